# Add "NA" values for duplicate_image_filename (column E) on rows 2 through 21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"

# Restore F1 (an empty-string cell in the original) which the runtime otherwise
# re-materializes with a stray value when shared strings are rewritten.
$ws.Range("F1").Value = ""
